$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of row 5 down into the new row 6 before we overwrite
# row 5's contents, so row 6 ends up with the same cell style (s="1").
$ws.Range("A5:C5").Copy()
$ws.Range("A6:C6").PasteSpecial(-4122)  # xlPasteFormats

# Turn row 1 into a header row: Name / Age / Technology
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Age"
$ws.Range("C1").Value = "Technology"

# Row 5 previously held "Malinikesh Agrawani" data; it now becomes
# Arindita Saha's row (with the corrected "Angular 2" spelling).
$ws.Range("A5").Value = "Arindita Saha"
$ws.Range("B5").Value = 21.0
$ws.Range("C5").Value = "Angular 2"

# New row 6 holds the "Malinikesh Agrawani" data that used to be on row 5,
# with the age corrected to 22.
$ws.Range("A6").Value = "Malinikesh Agrawani"
$ws.Range("B6").Value = 22.0
$ws.Range("C6").Value = "DCA"
